$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ODI Batting")
$new = $wb.Worksheets.Add($ws1)
$new.Name = "Test"
$new.Range("A1").Value = "x"
$new.Cells.Item(1,2) = ""
Write-Host "B1:" $new.Cells.Item(1,2).Value2
